# Update popup/parsing menu sheet with a "day number" (Cislo_dne) column,
# and roll the served dates forward by one week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H plus the per-row day index (1..5).
$ws.Range("H1").Value = "Cislo_dne"
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("H4").Value = 3
$ws.Range("H5").Value = 4
$ws.Range("H6").Value = 5

# Give the new column a sensible width (closest snap to the authored 10.22 chars).
$ws.Columns.Item(8).ColumnWidth = 9.3333333333333339

# Bump the "Od" (from) / "Do" (to) date range forward a week.
$ws.Range("A2").Value = 45767
$ws.Range("B2").Value = 45772

# Match the author's final selection in the saved file.
$ws.Range("L4").Select() | Out-Null
